$d = $word.ActiveDocument

$replacements = @(
    @("387÷3=", "454÷7="),
    @("880÷7=", "108÷3="),
    @("825÷3=", "858÷2="),
    @("505÷2=", "394÷9="),
    @("340÷5=", "975÷2="),
    @("516÷6=", "598÷3="),
    @("783÷5=", "544÷9="),
    @("218÷6=", "395÷9="),
    @("237÷2=", "919÷6="),
    @("855÷8=", "763÷3="),
    @("779÷8=", "188÷3="),
    @("850÷3=", "219÷3="),
    @("934÷5=", "990÷9="),
    @("966÷2=", "523÷3="),
    @("672÷3=", "631÷7="),
    @("931÷4=", "981÷3="),
    @("744÷8=", "455÷9="),
    @("782÷5=", "311÷6="),
    @("364÷4=", "555÷9="),
    @("848÷2=", "114÷5="),
    @("689÷8=", "775÷7="),
    @("810÷5=", "475÷3="),
    @("516÷4=", "857÷6="),
    @("566÷2=", "322÷8="),
    @("223÷2=", "980÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
